$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ajo" (Garlic) at
# "Feria Lagunitas de Puerto Montt". It is inserted as a new row 439,
# pushing the existing rows 439:530 down to 440:531 (dimension grows to
# A1:R531).
$ws.Rows.Item(439).Insert()

$ws.Range("A439").Value = 4
$ws.Range("B439").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C439").Value = "Los Lagos"
$ws.Range("D439").Value = Get-Date -Year 2023 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("E439").Value = 10
$ws.Range("F439").Value = 100112003
$ws.Range("G439").Value = "Ajo"
$ws.Range("H439").Value = "Chino"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 80
$ws.Range("K439").Value = 25000
$ws.Range("L439").Value = 26000
$ws.Range("M439").Value = 25500
$ws.Range("N439").Value = "$/caja 10 kilos"
$ws.Range("O439").Value = "China"
$ws.Range("P439").Value = 2550
$ws.Range("Q439").Value = 10
$ws.Range("R439").Value = "Hortaliza"
